# Applies the commit "need to figure out the stance and swing stats":
#  - Cuts the scratch AG:BE block out of the "alldata_1step" sheet
#  - Rebuilds/pastes a cleaned-up version of that block (split into a
#    "swing_df" table and a "stance_df" table) into the "pythonout" sheet
#    in columns K:O
#  - Leaves the workbook with "alldata_1step" as the active/selected sheet

$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("alldata_1step")
$ws6 = $wb.Worksheets.Item("pythonout")

# ------------------------------------------------------------------
# 1. Remove the old scratch data block (AG1:BE34) from alldata_1step
# ------------------------------------------------------------------
$ws4.Range("AG1:BE34").ClearContents()

# --- swing_df block (K41:O58) ---
$ws6.Range("K41").Value = "swing_df"
$ws6.Range("L42").Value = "subjectname"
$ws6.Range("M42").Value = "condname"
$ws6.Range("N42").Value = "trialname"
$ws6.Range("O42").Value = "metabolics_swing_avg_mean"
$ws6.Range("K43").Value = 0
$ws6.Range("L43").Value = "welk002"
$ws6.Range("M43").Value = "welkexo"
$ws6.Range("N43").Value = "trial01"
$ws6.Range("O43").Value = 2.8224860000000001
$ws6.Range("K44").Value = 1
$ws6.Range("L44").Value = "welk002"
$ws6.Range("M44").Value = "welkexo"
$ws6.Range("N44").Value = "trial02"
$ws6.Range("O44").Value = 2.7806150000000001
$ws6.Range("K45").Value = 2
$ws6.Range("L45").Value = "welk002"
$ws6.Range("M45").Value = "welkexo"
$ws6.Range("N45").Value = "trial03"
$ws6.Range("O45").Value = 2.5526460000000002
$ws6.Range("K46").Value = 3
$ws6.Range("L46").Value = "welk002"
$ws6.Range("M46").Value = "welkexo"
$ws6.Range("N46").Value = "trial04"
$ws6.Range("O46").Value = 2.8426930000000001
$ws6.Range("K47").Value = 4
$ws6.Range("L47").Value = "welk002"
$ws6.Range("M47").Value = "welknatural"
$ws6.Range("N47").Value = "trial01"
$ws6.Range("O47").Value = 2.670731
$ws6.Range("K48").Value = 5
$ws6.Range("L48").Value = "welk002"
$ws6.Range("M48").Value = "welknatural"
$ws6.Range("N48").Value = "trial02"
$ws6.Range("O48").Value = 3.1026220000000002
$ws6.Range("K49").Value = 6
$ws6.Range("L49").Value = "welk002"
$ws6.Range("M49").Value = "welknatural"
$ws6.Range("N49").Value = "trial03"
$ws6.Range("O49").Value = 2.6562060000000001
$ws6.Range("K50").Value = 7
$ws6.Range("L50").Value = "welk002"
$ws6.Range("M50").Value = "welknatural"
$ws6.Range("N50").Value = "trial04"
$ws6.Range("O50").Value = 2.774289
$ws6.Range("K51").Value = 8
$ws6.Range("L51").Value = "welk003"
$ws6.Range("M51").Value = "welkexo"
$ws6.Range("N51").Value = "trial01"
$ws6.Range("O51").Value = 3.190118
$ws6.Range("K52").Value = 9
$ws6.Range("L52").Value = "welk003"
$ws6.Range("M52").Value = "welkexo"
$ws6.Range("N52").Value = "trial02"
$ws6.Range("O52").Value = 3.1468889999999998
$ws6.Range("K53").Value = 10
$ws6.Range("L53").Value = "welk003"
$ws6.Range("M53").Value = "welkexo"
$ws6.Range("N53").Value = "trial03"
$ws6.Range("O53").Value = 3.2110289999999999
$ws6.Range("K54").Value = 11
$ws6.Range("L54").Value = "welk003"
$ws6.Range("M54").Value = "welkexo"
$ws6.Range("N54").Value = "trial04"
$ws6.Range("O54").Value = 3.5255000000000001
$ws6.Range("K55").Value = 12
$ws6.Range("L55").Value = "welk003"
$ws6.Range("M55").Value = "welknatural"
$ws6.Range("N55").Value = "trial01"
$ws6.Range("O55").Value = 3.3469850000000001
$ws6.Range("K56").Value = 13
$ws6.Range("L56").Value = "welk003"
$ws6.Range("M56").Value = "welknatural"
$ws6.Range("N56").Value = "trial02"
$ws6.Range("O56").Value = 3.0652469999999998
$ws6.Range("K57").Value = 14
$ws6.Range("L57").Value = "welk003"
$ws6.Range("M57").Value = "welknatural"
$ws6.Range("N57").Value = "trial03"
$ws6.Range("O57").Value = 3.0850390000000001
$ws6.Range("K58").Value = 15
$ws6.Range("L58").Value = "welk003"
$ws6.Range("M58").Value = "welknatural"
$ws6.Range("N58").Value = "trial04"
$ws6.Range("O58").Value = 3.1780020000000002

# --- stance_df block (K59:O76) ---
$ws6.Range("K59").Value = "stance_df"
$ws6.Range("L60").Value = "subjectname"
$ws6.Range("M60").Value = "condname"
$ws6.Range("N60").Value = "trialname"
$ws6.Range("O60").Value = "metabolics_stance_avg_mean"
$ws6.Range("K61").Value = 0
$ws6.Range("L61").Value = "welk002"
$ws6.Range("M61").Value = "welkexo"
$ws6.Range("N61").Value = "trial01"
$ws6.Range("O61").Value = 6.6376109999999997
$ws6.Range("K62").Value = 1
$ws6.Range("L62").Value = "welk002"
$ws6.Range("M62").Value = "welkexo"
$ws6.Range("N62").Value = "trial02"
$ws6.Range("O62").Value = 5.9639850000000001
$ws6.Range("K63").Value = 2
$ws6.Range("L63").Value = "welk002"
$ws6.Range("M63").Value = "welkexo"
$ws6.Range("N63").Value = "trial03"
$ws6.Range("O63").Value = 6.6552610000000003
$ws6.Range("K64").Value = 3
$ws6.Range("L64").Value = "welk002"
$ws6.Range("M64").Value = "welkexo"
$ws6.Range("N64").Value = "trial04"
$ws6.Range("O64").Value = 5.8253329999999997
$ws6.Range("K65").Value = 4
$ws6.Range("L65").Value = "welk002"
$ws6.Range("M65").Value = "welknatural"
$ws6.Range("N65").Value = "trial01"
$ws6.Range("O65").Value = 7.4643920000000001
$ws6.Range("K66").Value = 5
$ws6.Range("L66").Value = "welk002"
$ws6.Range("M66").Value = "welknatural"
$ws6.Range("N66").Value = "trial02"
$ws6.Range("O66").Value = 7.6236660000000001
$ws6.Range("K67").Value = 6
$ws6.Range("L67").Value = "welk002"
$ws6.Range("M67").Value = "welknatural"
$ws6.Range("N67").Value = "trial03"
$ws6.Range("O67").Value = 6.9179539999999999
$ws6.Range("K68").Value = 7
$ws6.Range("L68").Value = "welk002"
$ws6.Range("M68").Value = "welknatural"
$ws6.Range("N68").Value = "trial04"
$ws6.Range("O68").Value = 6.8632220000000004
$ws6.Range("K69").Value = 8
$ws6.Range("L69").Value = "welk003"
$ws6.Range("M69").Value = "welkexo"
$ws6.Range("N69").Value = "trial01"
$ws6.Range("O69").Value = 6.9033720000000001
$ws6.Range("K70").Value = 9
$ws6.Range("L70").Value = "welk003"
$ws6.Range("M70").Value = "welkexo"
$ws6.Range("N70").Value = "trial02"
$ws6.Range("O70").Value = 6.647322
$ws6.Range("K71").Value = 10
$ws6.Range("L71").Value = "welk003"
$ws6.Range("M71").Value = "welkexo"
$ws6.Range("N71").Value = "trial03"
$ws6.Range("O71").Value = 6.3015169999999996
$ws6.Range("K72").Value = 11
$ws6.Range("L72").Value = "welk003"
$ws6.Range("M72").Value = "welkexo"
$ws6.Range("N72").Value = "trial04"
$ws6.Range("O72").Value = 6.898752
$ws6.Range("K73").Value = 12
$ws6.Range("L73").Value = "welk003"
$ws6.Range("M73").Value = "welknatural"
$ws6.Range("N73").Value = "trial01"
$ws6.Range("O73").Value = 7.3225720000000001
$ws6.Range("K74").Value = 13
$ws6.Range("L74").Value = "welk003"
$ws6.Range("M74").Value = "welknatural"
$ws6.Range("N74").Value = "trial02"
$ws6.Range("O74").Value = 7.7138470000000003
$ws6.Range("K75").Value = 14
$ws6.Range("L75").Value = "welk003"
$ws6.Range("M75").Value = "welknatural"
$ws6.Range("N75").Value = "trial03"
$ws6.Range("O75").Value = 7.5333889999999997
$ws6.Range("K76").Value = 15
$ws6.Range("L76").Value = "welk003"
$ws6.Range("M76").Value = "welknatural"
$ws6.Range("N76").Value = "trial04"
$ws6.Range("O76").Value = 7.8034150000000002

# ------------------------------------------------------------------
# 2. Update the view state: select the new block on pythonout first,
#    then finish with alldata_1step active/selected (matches the
#    saved workbook view in the target).
# ------------------------------------------------------------------
$ws6.Select() | Out-Null
$ws6.Range("P43").Select() | Out-Null

$ws4.Select() | Out-Null
$ws4.Range("AG1:AQ16").Select() | Out-Null
